$wb = $excel.ActiveWorkbook

# zh-cn sheet: update Correspond Handoff/Handback DateTime for the
# 60262d6c-... row (row 3)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-18 14:35:25"
$wsZhCn.Range("H3").Value = "2016-03-18 14:35:58"

# de-de sheet: same row, same file, different language
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-18 14:35:33"
$wsDeDe.Range("H3").Value = "2016-03-18 14:36:04"
